# "added spark motor failure files"
#
# Semantic edit recovered from the OOXML diff (the rest of the diff -- Excel
# build/rupBuild, absolute author path, revision GUID, window geometry, and
# the row-height/column-bestFit "noise" coming from those -- are artifacts of
# the authoring machine/Excel build and are not deliberate worksheet edits):
#
#   1. Header cell C3 renamed "NN x KNN" -> "1-NN x K-NN"
#   2. The repeated sub-header label in column C (rows 5,7,9,11,13,15,17,19)
#      renamed "n" -> "k"
#   3. Active selection moved from I12 to A2
#   4. Columns re-sized (best effort; bestFit columns A:C narrower, D:K wider)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- 1. Header rename --------------------------------------------------
$ws.Range("C3").Value = "1-NN x K-NN"

# --- 2. Sub-header label rename (every block repeats this label) -------
$labelRows = @(5, 7, 9, 11, 13, 15, 17, 19)
foreach ($r in $labelRows) {
    $ws.Cells.Item($r, 3).Value = "k"
}

# --- 4. Column widths (best-effort match of the new layout) ------------
$ws.Columns.Item(1).ColumnWidth = 17.61
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 11.39
$ws.Range("D1:K1").EntireColumn.ColumnWidth = 16.94

# --- 3. Selection --------------------------------------------------------
[void]$ws.Range("A2").Select()
